$d = $word.ActiveDocument

# Original paragraph text is "Version 1." (indices: V0 e1 r2 s3 i4 o5 n6 _7 1_8 .9)
#
# Target XML splits the "Version" run into "Versi" + "on", changes
# "1" to "2", and moves the trailing "." into its own run placed
# after the _GoBack bookmark.

# Step 1: split "Version" into "Versi" | "on".
# Adding then immediately deleting a bookmark at the split point forces
# Word to break the run in two without leaving any residual formatting.
$splitPoint = $d.Range(5, 5)
$d.Bookmarks.Add("tmpSplit1", $splitPoint)
$d.Bookmarks("tmpSplit1").Delete()

# Step 2: change the "1" digit to "2"
$digit = $d.Range(8, 9)
$digit.Text = "2"

# Step 3: remove the trailing "." (currently at the end of the " 2." run)
$trailingDot = $d.Range(9, 10)
$trailingDot.Text = ""

# Step 4: re-insert "." as its own run after the _GoBack bookmark (end of content)
$tail = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$tail.InsertAfter(".")
